$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 9
$ws.Range("D2").Value = 44446
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 13000
$ws.Range("M2").Value = 13000
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 520

# Row 3 <- original row 6
$ws.Range("D3").Value = 44453
$ws.Range("J3").Value = 55
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14455
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 578

# Row 4 <- original row 5
$ws.Range("D4").Value = 44449
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 16000
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 640

# Row 5 <- original row 13
$ws.Range("D5").Value = 44421
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 600

# Row 6 <- original row 8
$ws.Range("D6").Value = 44432
$ws.Range("J6").Value = 15
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 14000
$ws.Range("O6").Value = "Provincia del Elquí"
$ws.Range("P6").Value = 560

# Row 7 <- original row 10
$ws.Range("D7").Value = 45113
$ws.Range("J7").Value = 8
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 17000
$ws.Range("M7").Value = 17000
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 680

# Row 8 <- original row 4
$ws.Range("D8").Value = 44340
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 600

# Row 9 <- original row 14
$ws.Range("D9").Value = 44376
$ws.Range("J9").Value = 15
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 480

# Row 10 <- original row 2
$ws.Range("D10").Value = 44467
$ws.Range("J10").Value = 35
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 480

# Row 11 <- original row 7
$ws.Range("D11").Value = 44418
$ws.Range("J11").Value = 12
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 600

# Row 12 <- original row 11
$ws.Range("D12").Value = 44435
$ws.Range("J12").Value = 15
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 14000
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 560

# Row 13 <- original row 12
$ws.Range("D13").Value = 44435
$ws.Range("J13").Value = 15
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 14000
$ws.Range("O13").Value = "Provincia del Elquí"
$ws.Range("P13").Value = 560

# Row 14 <- original row 3
$ws.Range("D14").Value = 44425
$ws.Range("J14").Value = 25
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 14000
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 560
